$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text cells (Coin name, Link, Volume label) ---
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("E10").Value = "9WazirXWRX"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"
$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("E13").Value = "12BitrueCoinBTR"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("E14").Value = "13BitMartTokenBMX"
$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("E15").Value = "14MCDexMCB"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("E16").Value = "15BitForexTokenBF"
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("E17").Value = "16CoinExTokenCET"
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("E20").Value = "19HotbitTokenHTBBestin24h"

# --- Numeric-looking text cells (Price column) must stay as text, matching the
#     original inlineStr cell type rather than being auto-converted to Number ---
$priceUpdates = @{
    "D2" = "243.03"
    "D4" = "5.391"
    "D5" = "0.05987"
    "D7" = "6.398"
    "D8" = "0.8115"
    "D9" = "0.8962"
    "D10" = "0.1410"
    "D11" = "0.07404"
    "D12" = "0.03379"
    "D13" = "0.03073"
    "D14" = "0.09330"
    "D15" = "3.868"
    "D16" = "0.001585"
    "D17" = "0.04649"
    "D18" = "0.0005939"
    "D19" = "0.006085"
    "D20" = "0.005016"
    "D21" = "0.0009840"
    "D22" = "0.00007798"
    "D23" = "0.0002899"
    "D24" = "3.614"
    "D25" = "2.161"
    "D27" = "0.1302"
    "D40" = "0.03893"
    "D41" = "0.006194"
    "D42" = "0.1075"
    "D43" = "0.002619"
    "D44" = "0.007187"
    "D45" = "0.00005190"
    "D47" = "0.0005799"
    "D48" = "0.9108"
    "D49" = "0.002298"
}
foreach ($cellRef in $priceUpdates.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$cellRef]
    $cell.Style = "Normal"
}

Write-Host "Applied cryptos.xlsx price/listing updates"
